$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.847.57"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "3.404.56"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.12"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.68"
$ws.Range("E6").Value = "  -3.61%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.405.09"
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.57"
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.125"
$ws.Range("E11").Value = "  +1.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.395"
$ws.Range("E12").Value = "  +2.55%  "
$ws.Range("D13").Value = "3.981.11"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.26"
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000173"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "3.417.12"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "60.959.19"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.32"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.27"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.25"
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "389.63"
$ws.Range("E22").Value = "  +1.85%  "
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.94"
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000123"
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").Value = "3.548.97"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.180"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.40"
$ws.Range("E30").Value = "  -4.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.18"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("E32").Value = "  -6.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.17"
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.85"
$ws.Range("E35").Value = "  -0.52%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").Value = "3.430.18"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "167.38"
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.54"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0785"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.90"
$ws.Range("E42").Value = "  +4.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.790"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.49"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.89"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.70"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("D48").Value = "2.579.82"
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.13"
$ws.Range("E49").Value = "  -3.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.93"
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.08"
$ws.Range("E51").Value = "  -1.96%  "
